# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, exactly as the source feed renders it
# (Price/Volume columns are free-text strings, not numeric cells).
$updates = [ordered]@{
    'D2' = '66.161.70'
    'E2' = '  +7.22%  '
    'D3' = '3.017.90'
    'E3' = '  +4.38%  '
    'E4' = '  +0.10%  '
    'D5' = '585.10'
    'E5' = '  +3.49%  '
    'D6' = '156.18'
    'E6' = '  +9.33%  '
    'E7' = '  -0.02%  '
    'D8' = '3.011.03'
    'E8' = '  +4.14%  '
    'E9' = '  +3.53%  '
    'E10' = '  +1.06%  '
    'E11' = '  +6.94%  '
    'D12' = '0.452'
    'E12' = '  +5.87%  '
    'D13' = '0.0000251'
    'E13' = '  +9.01%  '
    'D14' = '34.65'
    'E14' = '  +9.60%  '
    'E15' = '  +0.59%  '
    'D16' = '66.176.47'
    'E16' = '  +7.33%  '
    'D17' = '3.516.56'
    'E17' = '  +4.31%  '
    'D18' = '6.96'
    'E18' = '  +7.07%  '
    'D19' = '3.011.54'
    'E19' = '  +3.86%  '
    'D20' = '462.66'
    'E20' = '  +7.83%  '
    'D21' = '13.88'
    'E21' = '  +6.82%  '
    'E22' = '  +5.41%  '
    'E23' = '  +8.59%  '
    'D24' = '82.11'
    'E24' = '  +4.05%  '
    'E25' = '  +13.34%  '
    'D26' = '12.49'
    'E26' = '  +5.07%  '
    'D27' = '10.71'
    'E27' = '  +8.23%  '
    'E28' = '  -0.10%  '
    'D29' = '8.03'
    'E29' = '  +14.94%  '
    'E30' = '  +16.14%  '
    'D31' = '0.0000105'
    'E31' = '  -1.28%  '
    'E32' = '  +5.39%  '
    'D33' = '0.111'
    'E33' = '  +5.15%  '
    'D34' = '27.02'
    'E34' = '  +6.28%  '
    'D35' = '1.00'
    'E35' = '  -0.04%  '
    'D36' = '0.996'
    'E36' = '  +4.14%  '
    'D37' = '5.82'
    'E37' = '  +8.75%  '
    'E38' = '  +12.93%  '
    'E39' = '  +10.14%  '
    'D40' = '49.52'
    'E40' = '  +1.50%  '
    'E41' = '  +8.66%  '
    'E42' = '  +14.01%  '
    'D43' = '43.70'
    'E43' = '  +11.25%  '
    'E44' = '  +3.60%  '
    'D45' = '393.97'
    'E45' = '  +14.96%  '
    'D46' = '2.801.00'
    'E46' = '  +4.52%  '
    'E47' = '  +5.95%  '
    'D48' = '133.92'
    'E48' = '  +1.34%  '
    'D50' = '23.58'
    'E50' = '  +9.76%  '
    'E51' = '  +4.33%  '
}

foreach ($ref in $updates.Keys) {
    $newValue = $updates[$ref]
    $cell = $ws.Range($ref)

    # Numeric-looking text (e.g. '1.00', '49.52') must stay text so Excel doesn't
    # silently coerce it to a Double and drop formatting like trailing zeros.
    if ($newValue.Trim() -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = '@'
        $cell.Value = $newValue
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $newValue
    }
}
